# Applies the "Add files via upload" commit:
#  - inserts a new "ProviderOptions" worksheet (sheetId 10) right before
#    "PlacementPreservationStrategy"
#  - populates its two template rows (header row + the two standard
#    testPOM / testT4273 rows used throughout this workbook)
#  - nudges a handful of cached cell-selections / fills a couple of blank
#    template cells on other sheets, and leaves ChildLocationCorrespondences
#    as the active tab (matching the final saved UI state captured in the diff)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "ProviderOptions" sheet immediately before
#    "PlacementPreservationStrategy" (this lands it at position 5, gets
#    sheetId 10, and cascades the rId/physical sheetN.xml renumbering for
#    every sheet after it).
# ---------------------------------------------------------------------
$placementPreservation = $wb.Worksheets.Item("PlacementPreservationStrategy")
$providerOptions = $wb.Worksheets.Add($placementPreservation)
$providerOptions.Name = "ProviderOptions"

# Reference sheets whose existing cells already carry the two cell styles
# ("header" style s=1, and the gray "record row" style s=3) we need to
# reuse on the new sheet, so we copy the formatting instead of trying to
# rebuild it from scratch.
$removal = $wb.Worksheets.Item("Removal")

# Header row formatting (style s=1) across A1:S1
$removal.Range("A1").Copy()
$providerOptions.Range("A1:S1").PasteSpecial(-4122)

# Header values (H1 intentionally stays blank, matching the template)
$providerOptions.Range("A1").Value = "TESTCASE"
$providerOptions.Range("B1").Value = "SCRIPT_ITERATION"
$providerOptions.Range("C1").Value = "POM_ITERATION"
$providerOptions.Range("D1").Value = "PLACEMENT_SERVICE_PROVIDER_OPTIONS"
$providerOptions.Range("E1").Value = "PROVIDER_OPTIONS"
$providerOptions.Range("F1").Value = "SERVICE_OPTIONS"
$providerOptions.Range("G1").Value = "PROVIDER_OPTION_DETERMINATION_STATUS"
$providerOptions.Range("I1").Value = "SAVE_BTN"
$providerOptions.Range("J1").Value = "NEW_ITEM_SERVICE_REQUEST"
$providerOptions.Range("K1").Value = "STATUS"
$providerOptions.Range("L1").Value = "APPROVAL_SUPERVISOR"
$providerOptions.Range("M1").Value = "APPROVAL_SUPERVISOR_LINK"
$providerOptions.Range("N1").Value = "PLACEMENT_OPTIONS_RECORD"
$providerOptions.Range("O1").Value = "EMERGENCY_CHECKBOX"
$providerOptions.Range("P1").Value = "PLACEMENT_OPTION_DETERMINATION_STATUS"
$providerOptions.Range("Q1").Value = "WILL_PLACEMENT_PROCEED"
$providerOptions.Range("R1").Value = "REASON"
$providerOptions.Range("S1").Value = "OTHER_REASON"

# Row 2 ("testPOM" smoke-test row) - plain, unstyled cells
$providerOptions.Range("A2").Value = "testPOM"
$providerOptions.Range("B2").Value = 1
$providerOptions.Range("C2").Value = 1

# Row 3 ("testT4273" row) - A3 uses the gray record-row style (s=3);
# B3/C3 stay plain, matching the target workbook exactly.
$removal.Range("A3").Copy()
$providerOptions.Range("A3").PasteSpecial(-4122)
$providerOptions.Range("A3").Value = "testT4273"
$providerOptions.Range("B3").Value = 1
$providerOptions.Range("C3").Value = 1

# Leave this sheet's own cached selection on H1 (the blank header cell)
$providerOptions.Range("H1").Select()

# ---------------------------------------------------------------------
# 2. FolioChildLocationAbsences: cached selection moves from "select all"
#    to a single cell, D9.
# ---------------------------------------------------------------------
$absences = $wb.Worksheets.Item("FolioChildLocationAbsences")
$absences.Activate()
$absences.Range("D9").Select()

# ---------------------------------------------------------------------
# 3. PlacementPreservationStrategy: cached selection moves from
#    "select all" to A1:C2.
# ---------------------------------------------------------------------
$placementPreservation.Activate()
$placementPreservation.Range("A1:C2").Select()

# ---------------------------------------------------------------------
# 4. OverstayReportDetails: cached selection moves from A1:D3 to A3.
# ---------------------------------------------------------------------
$overstayDetails = $wb.Worksheets.Item("OverstayReportDetails")
$overstayDetails.Activate()
$overstayDetails.Range("A3").Select()

# ---------------------------------------------------------------------
# 5. ChildLocationCorrespondences: fill in the two blank template cells
#    on the "testT4273" row (J3/K3), matching the same values already
#    used on the row above (row 2), then leave this as the active tab
#    with K3 selected - matching the final saved workbook state.
# ---------------------------------------------------------------------
$correspondences = $wb.Worksheets.Item("ChildLocationCorrespondences")
$correspondences.Activate()
$correspondences.Range("J3").Value = "past"
$correspondences.Range("K3").Value = "Court Ordered Placement"
$correspondences.Range("K3").Select()
